# Add a new daily-report sheet "1.11" by duplicating the "1.10)" sheet
# (which holds the previous day's report) and updating it with the new
# day's data, mirroring the author's "Add files via upload" commit.

$wb = $excel.ActiveWorkbook

# The existing "1.10)" sheet is the template/source for the new day.
$source = $wb.Worksheets.Item("1.10)")
$source.Activate()

# Before duplicating, restore the source sheet's on-screen selection to
# the range it ends up showing once it is no longer the active tab.
$source.Range("C10:E10").Select()

# Duplicate "1.10)" and place the copy right after it; Excel names the
# copy "1.10) (2)" and makes it the new active sheet automatically.
$source.Copy($null, $source)
$newSheet = $wb.Worksheets.Item("1.10) (2)")
$newSheet.Name = "1.11"

# Update the new sheet's header date (format D3 "日期：" / E3 value) and
# the planned/actual date rows from "1.1" to "1.11".
$newSheet.Range("E3").Value = 1.11
$newSheet.Range("C6").Value = 1.11
$newSheet.Range("C7").Value = 1.11

# Record today's (1.11) work log entries, replacing yesterday's text.
$newSheet.Range("G10").Value = "学习了QtAction"
$newSheet.Range("G11").Value = "学习了歌词的显示"

# Task 2 was not completed today.
$newSheet.Range("J11").Value = "N"

# Keep the new sheet's selection where the author last left the cursor.
$newSheet.Range("J11").Select()
